{"js": "// Append the ASP.NET/HTML \"Example.aspx\" markup as a new block of plain-text\n// paragraphs after the existing pageLoad() code sample (mirrors the\n// UpdatePanel/EndRequestHandler sample added in the source diff).\nconst newParagraphs = [\n  \"<%@ Page Language=\\\"C#\\\" AutoEventWireup=\\\"true\\\" CodeBehind=\\\"Example.aspx.cs\\\" Inherits=\\\"YourNamespace.Example\\\" %>\",\n  \"\",\n  \"<!DOCTYPE html>\",\n  \"\",\n  \"<html xmlns=\\\"http://www.w3.org/1999/xhtml\\\">\",\n  \"<head runat=\\\"server\\\">\",\n  \"    <title></title>\",\n  \"    <script src=\\\"https://ajax.googleapis.com/ajax/libs/jquery/3.5.1/jquery.min.js\\\"></script>\",\n  \"    <script type=\\\"text/javascript\\\">\",\n  \"        Sys.WebForms.PageRequestManager.getInstance().add_endRequest(EndRequestHandler);\",\n  \"\",\n  \"        function EndRequestHandler(sender, args) {\",\n  \"            // This function will run after every partial postback\",\n  \"            $('#divToUpdate').attr('customAttribute', 'customValue');\",\n  \"        }\",\n  \"    </script>\",\n  \"</head>\",\n  \"<body>\",\n  \"    <form id=\\\"form1\\\" runat=\\\"server\\\">\",\n  \"        <asp:ScriptManager ID=\\\"ScriptManager1\\\" runat=\\\"server\\\"></asp:ScriptManager>\",\n  \"        <asp:UpdatePanel ID=\\\"UpdatePanel1\\\" runat=\\\"server\\\">\",\n  \"            <ContentTemplate>\",\n  \"                <div id=\\\"divToUpdate\\\">\",\n  \"                    This is a div inside the UpdatePanel.\",\n  \"                </div>\",\n  \"                <asp:Button ID=\\\"btnUpdate\\\" runat=\\\"server\\\" Text=\\\"Update\\\" OnClick=\\\"btnUpdate_Click\\\" />\",\n  \"            </ContentTemplate>\",\n  \"        </asp:UpdatePanel>\",\n  \"    </form>\",\n  \"</body>\",\n  \"</html>\"\n];\n\n// Find the last paragraph in the document body (\"}\" that closes pageLoad)\n// and insert the new paragraphs, one by one, right after it (and therefore\n// before the section properties / end of body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet last = paragraphs.items[paragraphs.items.length - 1];\nfor (const text of newParagraphs) {\n  last = last.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Append the ASP.NET/HTML \"Example.aspx\" markup as a new block of plain-text\n# paragraphs after the existing pageLoad() code sample (mirrors the\n# UpdatePanel/EndRequestHandler sample added in the source diff).\n$d = $word.ActiveDocument\n\n$newParagraphs = @(\n    '<%@ Page Language=\"C#\" AutoEventWireup=\"true\" CodeBehind=\"Example.aspx.cs\" Inherits=\"YourNamespace.Example\" %>',\n    '',\n    '<!DOCTYPE html>',\n    '',\n    '<html xmlns=\"http://www.w3.org/1999/xhtml\">',\n    '<head runat=\"server\">',\n    '    <title></title>',\n    '    <script src=\"https://ajax.googleapis.com/ajax/libs/jquery/3.5.1/jquery.min.js\"></script>',\n    '    <script type=\"text/javascript\">',\n    '        Sys.WebForms.PageRequestManager.getInstance().add_endRequest(EndRequestHandler);',\n    '',\n    '        function EndRequestHandler(sender, args) {',\n    '            // This function will run after every partial postback',\n    '            $(''#divToUpdate'').attr(''customAttribute'', ''customValue'');',\n    '        }',\n    '    </script>',\n    '</head>',\n    '<body>',\n    '    <form id=\"form1\" runat=\"server\">',\n    '        <asp:ScriptManager ID=\"ScriptManager1\" runat=\"server\"></asp:ScriptManager>',\n    '        <asp:UpdatePanel ID=\"UpdatePanel1\" runat=\"server\">',\n    '            <ContentTemplate>',\n    '                <div id=\"divToUpdate\">',\n    '                    This is a div inside the UpdatePanel.',\n    '                </div>',\n    '                <asp:Button ID=\"btnUpdate\" runat=\"server\" Text=\"Update\" OnClick=\"btnUpdate_Click\" />',\n    '            </ContentTemplate>',\n    '        </asp:UpdatePanel>',\n    '    </form>',\n    '</body>',\n    '</html>'\n)\n\n# Start from the last paragraph in the document (\"}\" that closes pageLoad)\n# and insert the new paragraphs, one by one, right after it (and therefore\n# before the section properties / end of body).\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$r = $lastPara.Range\n\nforeach ($text in $newParagraphs) {\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    $p = $d.Paragraphs.Item($d.Paragraphs.Count)\n    if ($text -ne '') {\n        $p.Range.Text = $text\n    }\n    $r = $p.Range\n}\n"}
